$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 136
$ws.Range("D136").Value = 44719
$ws.Range("I136").Value = '1a (guarda)'
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 620
$ws.Range("L136").Value = 650
$ws.Range("M136").Value = 633
$ws.Range("P136").Value = 633

# Row 137
$ws.Range("D137").Value = 44719
$ws.Range("H137").Value = 'Camote'
$ws.Range("I137").Value = '2a (guarda)'
$ws.Range("J137").Value = 400
$ws.Range("K137").Value = 500
$ws.Range("L137").Value = 500
$ws.Range("M137").Value = 500
$ws.Range("P137").Value = 500

# Row 138
$ws.Range("D138").Value = 44637
$ws.Range("J138").Value = 300
$ws.Range("M138").Value = 325
$ws.Range("O138").Value = 'Región de O''Higgins'
$ws.Range("P138").Value = 325

# Row 139
$ws.Range("D139").Value = 44461
$ws.Range("H139").Value = 'Paine'
$ws.Range("I139").Value = '1a (guarda)'
$ws.Range("J139").Value = 300
$ws.Range("K139").Value = 200
$ws.Range("L139").Value = 220
$ws.Range("M139").Value = 210
$ws.Range("O139").Value = 'Región de O''Higgins'
$ws.Range("P139").Value = 210

# Row 140
$ws.Range("D140").Value = 44600
$ws.Range("I140").Value = '1a (cosecha)'
$ws.Range("J140").Value = 500
$ws.Range("K140").Value = 300
$ws.Range("L140").Value = 350
$ws.Range("M140").Value = 320
$ws.Range("O140").Value = 'Región del Maule'
$ws.Range("P140").Value = 320

# Row 141
$ws.Range("D141").Value = 44629
$ws.Range("I141").Value = '1a (cosecha)'
$ws.Range("J141").Value = 500
$ws.Range("K141").Value = 300
$ws.Range("L141").Value = 350
$ws.Range("M141").Value = 330
$ws.Range("O141").Value = 'Región del Maule'
$ws.Range("P141").Value = 330

# Row 142
$ws.Range("D142").Value = 44364
$ws.Range("J142").Value = 600
$ws.Range("L142").Value = 370
$ws.Range("M142").Value = 360
$ws.Range("P142").Value = 360

# Row 143
$ws.Range("D143").Value = 44364
$ws.Range("K143").Value = 250
$ws.Range("L143").Value = 250
$ws.Range("M143").Value = 250
$ws.Range("P143").Value = 250

# Row 144
$ws.Range("D144").Value = 44341
$ws.Range("I144").Value = '1a (guarda)'
$ws.Range("J144").Value = 500
$ws.Range("K144").Value = 350
$ws.Range("L144").Value = 380
$ws.Range("M144").Value = 368
$ws.Range("P144").Value = 368

# Row 145
$ws.Range("D145").Value = 44341
$ws.Range("I145").Value = '2a (guarda)'
$ws.Range("J145").Value = 300
$ws.Range("K145").Value = 300
$ws.Range("L145").Value = 300
$ws.Range("M145").Value = 300
$ws.Range("P145").Value = 300

# Row 146
$ws.Range("D146").Value = 44649
$ws.Range("I146").Value = '1a (cosecha)'
$ws.Range("J146").Value = 700
$ws.Range("K146").Value = 300
$ws.Range("M146").Value = 325
$ws.Range("P146").Value = 325

# Row 147
$ws.Range("D147").Value = 44230
$ws.Range("I147").Value = '1a nueva(o)'
$ws.Range("K147").Value = 450
$ws.Range("L147").Value = 450
$ws.Range("M147").Value = 450
$ws.Range("P147").Value = 450

# Row 148
$ws.Range("D148").Value = 44230
$ws.Range("I148").Value = '2a nueva(o)'
$ws.Range("K148").Value = 350
$ws.Range("L148").Value = 350
$ws.Range("M148").Value = 350
$ws.Range("P148").Value = 350

# Row 149
$ws.Range("D149").Value = 44320
$ws.Range("I149").Value = '1a (cosecha)'
$ws.Range("J149").Value = 400
$ws.Range("K149").Value = 350
$ws.Range("L149").Value = 350
$ws.Range("M149").Value = 350
$ws.Range("P149").Value = 350

# Row 150
$ws.Range("D150").Value = 44320
$ws.Range("I150").Value = '2a (cosecha)'
$ws.Range("J150").Value = 300
$ws.Range("K150").Value = 250
$ws.Range("L150").Value = 250
$ws.Range("M150").Value = 250
$ws.Range("P150").Value = 250

# Row 151
$ws.Range("D151").Value = 44257
$ws.Range("I151").Value = '1a nueva(o)'
$ws.Range("J151").Value = 300
$ws.Range("K151").Value = 400
$ws.Range("L151").Value = 400
$ws.Range("M151").Value = 400
$ws.Range("P151").Value = 400

# Row 152
$ws.Range("D152").Value = 44257
$ws.Range("I152").Value = '2a nueva(o)'
$ws.Range("J152").Value = 200
$ws.Range("K152").Value = 320
$ws.Range("L152").Value = 320
$ws.Range("M152").Value = 320
$ws.Range("P152").Value = 320

# Row 153
$ws.Range("D153").Value = 44428
$ws.Range("I153").Value = '1a (guarda)'
$ws.Range("J153").Value = 600
$ws.Range("K153").Value = 600
$ws.Range("L153").Value = 650
$ws.Range("M153").Value = 625
$ws.Range("P153").Value = 625

# Row 154
$ws.Range("D154").Value = 44428
$ws.Range("I154").Value = '2a (guarda)'
$ws.Range("J154").Value = 300
$ws.Range("K154").Value = 500
$ws.Range("L154").Value = 500
$ws.Range("M154").Value = 500
$ws.Range("P154").Value = 500

# Row 155
$ws.Range("D155").Value = 44596
$ws.Range("I155").Value = '1a nueva(o)'
$ws.Range("J155").Value = 500
$ws.Range("K155").Value = 450
$ws.Range("L155").Value = 450
$ws.Range("M155").Value = 450
$ws.Range("O155").Value = 'Región de O''Higgins'
$ws.Range("P155").Value = 450

# Row 156
$ws.Range("D156").Value = 44596
$ws.Range("I156").Value = '2a nueva(o)'
$ws.Range("J156").Value = 500
$ws.Range("K156").Value = 350
$ws.Range("L156").Value = 350
$ws.Range("M156").Value = 350
$ws.Range("O156").Value = 'Región de O''Higgins'
$ws.Range("P156").Value = 350

# Row 157
$ws.Range("D157").Value = 44447
$ws.Range("J157").Value = 600
$ws.Range("K157").Value = 600
$ws.Range("L157").Value = 650
$ws.Range("M157").Value = 625
$ws.Range("O157").Value = 'Provincia de Melipilla'
$ws.Range("P157").Value = 625

# Row 158
$ws.Range("D158").Value = 44447
$ws.Range("K158").Value = 550
$ws.Range("L158").Value = 550
$ws.Range("M158").Value = 550
$ws.Range("O158").Value = 'Provincia de Melipilla'
$ws.Range("P158").Value = 550

# Row 159
$ws.Range("D159").Value = 44420
$ws.Range("J159").Value = 300
$ws.Range("K159").Value = 450
$ws.Range("L159").Value = 450
$ws.Range("M159").Value = 450
$ws.Range("P159").Value = 450

# Row 160
$ws.Range("D160").Value = 44420
$ws.Range("J160").Value = 300
$ws.Range("K160").Value = 350
$ws.Range("L160").Value = 350
$ws.Range("M160").Value = 350
$ws.Range("P160").Value = 350

# Row 161
$ws.Range("D161").Value = 44435
$ws.Range("I161").Value = '1a (guarda)'
$ws.Range("J161").Value = 500
$ws.Range("K161").Value = 850
$ws.Range("L161").Value = 900
$ws.Range("M161").Value = 880
$ws.Range("P161").Value = 880

# Row 162
$ws.Range("D162").Value = 44435
$ws.Range("I162").Value = '2a (guarda)'
$ws.Range("J162").Value = 200
$ws.Range("K162").Value = 800
$ws.Range("L162").Value = 800
$ws.Range("M162").Value = 800
$ws.Range("P162").Value = 800

# Row 163
$ws.Range("D163").Value = 44208
$ws.Range("J163").Value = 400
$ws.Range("L163").Value = 500
$ws.Range("M163").Value = 500
$ws.Range("N163").Value = '$/kilo (volumen en unidades)'
$ws.Range("P163").Value = 500

# Row 164
$ws.Range("D164").Value = 44208
$ws.Range("I164").Value = '2a nueva(o)'
$ws.Range("J164").Value = 300
$ws.Range("K164").Value = 400
$ws.Range("L164").Value = 400
$ws.Range("M164").Value = 400
$ws.Range("P164").Value = 400

# Row 165
$ws.Range("D165").Value = 44496
$ws.Range("I165").Value = '1a nueva(o)'
$ws.Range("J165").Value = 1000
$ws.Range("K165").Value = 500
$ws.Range("M165").Value = 560
$ws.Range("N165").Value = '$/kilo'
$ws.Range("P165").Value = 560

# Row 166
$ws.Range("D166").Value = 44449
$ws.Range("I166").Value = '1a (guarda)'
$ws.Range("K166").Value = 650
$ws.Range("L166").Value = 700
$ws.Range("M166").Value = 675
$ws.Range("P166").Value = 675

# Row 167
$ws.Range("D167").Value = 44449
$ws.Range("I167").Value = '2a (guarda)'
$ws.Range("J167").Value = 200
$ws.Range("K167").Value = 600
$ws.Range("L167").Value = 600
$ws.Range("M167").Value = 600
$ws.Range("P167").Value = 600

# Row 168
$ws.Range("D168").Value = 44204
$ws.Range("I168").Value = '1a nueva(o)'
$ws.Range("J168").Value = 400
$ws.Range("K168").Value = 450
$ws.Range("L168").Value = 450
$ws.Range("M168").Value = 450
$ws.Range("P168").Value = 450

# Row 169
$ws.Range("D169").Value = 44204
$ws.Range("I169").Value = '2a nueva(o)'
$ws.Range("J169").Value = 400
$ws.Range("K169").Value = 350
$ws.Range("L169").Value = 350
$ws.Range("M169").Value = 350
$ws.Range("P169").Value = 350

# Row 170
$ws.Range("D170").Value = 44679
$ws.Range("J170").Value = 300
$ws.Range("K170").Value = 650
$ws.Range("L170").Value = 650
$ws.Range("M170").Value = 650
$ws.Range("P170").Value = 650

# Row 171
$ws.Range("D171").Value = 44679
$ws.Range("K171").Value = 550
$ws.Range("L171").Value = 550
$ws.Range("M171").Value = 550
$ws.Range("P171").Value = 550

# Row 172
$ws.Range("D172").Value = 44663
$ws.Range("J172").Value = 600
$ws.Range("K172").Value = 350
$ws.Range("M172").Value = 375
$ws.Range("P172").Value = 375

# Row 173
$ws.Range("D173").Value = 44663

# Row 174
$ws.Range("D174").Value = 44313
$ws.Range("I174").Value = '1a (cosecha)'
$ws.Range("K174").Value = 400
$ws.Range("L174").Value = 400
$ws.Range("M174").Value = 400
$ws.Range("P174").Value = 400

# Row 175
$ws.Range("D175").Value = 44313
$ws.Range("I175").Value = '2a (cosecha)'
$ws.Range("K175").Value = 300
$ws.Range("L175").Value = 300
$ws.Range("M175").Value = 300
$ws.Range("P175").Value = 300

# Row 176
$ws.Range("D176").Value = 44195
$ws.Range("I176").Value = '1a nueva(o)'
$ws.Range("J176").Value = 300
$ws.Range("K176").Value = 800
$ws.Range("L176").Value = 800
$ws.Range("M176").Value = 800
$ws.Range("P176").Value = 800

# Row 177
$ws.Range("D177").Value = 44195
$ws.Range("I177").Value = '2a nueva(o)'
$ws.Range("J177").Value = 300
$ws.Range("K177").Value = 700
$ws.Range("L177").Value = 700
$ws.Range("M177").Value = 700
$ws.Range("P177").Value = 700

# Row 178
$ws.Range("D178").Value = 44433
$ws.Range("I178").Value = '1a (guarda)'
$ws.Range("J178").Value = 500
$ws.Range("K178").Value = 850
$ws.Range("L178").Value = 900
$ws.Range("M178").Value = 880
$ws.Range("O178").Value = 'Región de O''Higgins'
$ws.Range("P178").Value = 880

# Row 179
$ws.Range("D179").Value = 44433
$ws.Range("I179").Value = '2a (guarda)'
$ws.Range("K179").Value = 800
$ws.Range("L179").Value = 800
$ws.Range("M179").Value = 800
$ws.Range("O179").Value = 'Región de O''Higgins'
$ws.Range("P179").Value = 800

# Row 180
$ws.Range("D180").Value = 44292
$ws.Range("I180").Value = '1a (cosecha)'
$ws.Range("J180").Value = 250
$ws.Range("K180").Value = 400
$ws.Range("L180").Value = 400
$ws.Range("M180").Value = 400
$ws.Range("O180").Value = 'Región Metropolitana'
$ws.Range("P180").Value = 400

# Row 181
$ws.Range("D181").Value = 44292
$ws.Range("I181").Value = '2a (cosecha)'
$ws.Range("J181").Value = 200
$ws.Range("K181").Value = 300
$ws.Range("L181").Value = 300
$ws.Range("M181").Value = 300
$ws.Range("O181").Value = 'Región Metropolitana'
$ws.Range("P181").Value = 300

# Row 182
$ws.Range("D182").Value = 44579
$ws.Range("I182").Value = '1a nueva(o)'
$ws.Range("J182").Value = 800
$ws.Range("K182").Value = 500
$ws.Range("L182").Value = 550
$ws.Range("M182").Value = 525
$ws.Range("P182").Value = 525

# Row 183
$ws.Range("D183").Value = 44579
$ws.Range("I183").Value = '2a nueva(o)'
$ws.Range("J183").Value = 400
$ws.Range("K183").Value = 450
$ws.Range("L183").Value = 450
$ws.Range("M183").Value = 450
$ws.Range("P183").Value = 450

# Row 184
$ws.Range("H184").Value = 'Camote'
$ws.Range("K184").Value = 400
$ws.Range("L184").Value = 400
$ws.Range("M184").Value = 400
$ws.Range("P184").Value = 400

# Row 185
$ws.Range("H185").Value = 'Camote'
$ws.Range("K185").Value = 300
$ws.Range("L185").Value = 300
$ws.Range("M185").Value = 300
$ws.Range("P185").Value = 300

# Row 186
$ws.Range("D186").Value = 44285
$ws.Range("I186").Value = '1a (cosecha)'
$ws.Range("J186").Value = 300
$ws.Range("K186").Value = 250
$ws.Range("L186").Value = 250
$ws.Range("M186").Value = 250
$ws.Range("P186").Value = 250

# Row 187
$ws.Range("D187").Value = 44285
$ws.Range("H187").Value = 'Paine'
$ws.Range("I187").Value = '2a (cosecha)'
$ws.Range("K187").Value = 200
$ws.Range("L187").Value = 200
$ws.Range("M187").Value = 200
$ws.Range("O187").Value = 'Región de O''Higgins'
$ws.Range("P187").Value = 200

# Row 188
$ws.Range("D188").Value = 44160
$ws.Range("H188").Value = 'Paine'
$ws.Range("I188").Value = '1a nueva(o)'
$ws.Range("J188").Value = 500
$ws.Range("K188").Value = 800
$ws.Range("L188").Value = 1000
$ws.Range("M188").Value = 920
$ws.Range("O188").Value = 'Región de O''Higgins'
$ws.Range("P188").Value = 920

# Row 189
$ws.Range("D189").Value = 44554
$ws.Range("J189").Value = 300
$ws.Range("K189").Value = 650
$ws.Range("L189").Value = 650
$ws.Range("M189").Value = 650
$ws.Range("O189").Value = 'Región de Coquimbo'
$ws.Range("P189").Value = 650

# Row 190
$ws.Range("D190").Value = 44554
$ws.Range("J190").Value = 300
$ws.Range("K190").Value = 550
$ws.Range("L190").Value = 550
$ws.Range("M190").Value = 550
$ws.Range("O190").Value = 'Región de Coquimbo'
$ws.Range("P190").Value = 550

# Row 191
$ws.Range("D191").Value = 44565
$ws.Range("I191").Value = '1a nueva(o)'
$ws.Range("J191").Value = 800
$ws.Range("K191").Value = 500
$ws.Range("L191").Value = 550
$ws.Range("M191").Value = 525
$ws.Range("O191").Value = 'Región Metropolitana'
$ws.Range("P191").Value = 525

# Row 192
$ws.Range("D192").Value = 44565
$ws.Range("I192").Value = '2a nueva(o)'
$ws.Range("J192").Value = 400
$ws.Range("K192").Value = 450
$ws.Range("L192").Value = 450
$ws.Range("M192").Value = 450
$ws.Range("O192").Value = 'Región Metropolitana'
$ws.Range("P192").Value = 450

# Row 193
$ws.Range("D193").Value = 44603
$ws.Range("I193").Value = '1a (cosecha)'
$ws.Range("J193").Value = 350
$ws.Range("K193").Value = 300
$ws.Range("L193").Value = 350
$ws.Range("M193").Value = 321
$ws.Range("P193").Value = 321

# Row 194
$ws.Range("D194").Value = 44187
$ws.Range("I194").Value = '1a nueva(o)'
$ws.Range("J194").Value = 300
$ws.Range("K194").Value = 600
$ws.Range("L194").Value = 600
$ws.Range("M194").Value = 600
$ws.Range("O194").Value = 'Región de O''Higgins'
$ws.Range("P194").Value = 600

# Row 195
$ws.Range("D195").Value = 44187
$ws.Range("I195").Value = '2a nueva(o)'
$ws.Range("J195").Value = 200
$ws.Range("K195").Value = 500
$ws.Range("L195").Value = 500
$ws.Range("M195").Value = 500
$ws.Range("P195").Value = 500

# Row 196
$ws.Range("D196").Value = 44638
$ws.Range("I196").Value = '1a (cosecha)'
$ws.Range("J196").Value = 350
$ws.Range("K196").Value = 300
$ws.Range("L196").Value = 350
$ws.Range("M196").Value = 321
$ws.Range("O196").Value = 'Región del Maule'
$ws.Range("P196").Value = 321

# Row 197
$ws.Range("D197").Value = 44299
$ws.Range("J197").Value = 300
$ws.Range("L197").Value = 300
$ws.Range("M197").Value = 300
$ws.Range("P197").Value = 300

# Row 198
$ws.Range("D198").Value = 44299
$ws.Range("I198").Value = '2a (cosecha)'
$ws.Range("J198").Value = 300
$ws.Range("K198").Value = 250
$ws.Range("L198").Value = 250
$ws.Range("M198").Value = 250
$ws.Range("P198").Value = 250

# Row 199
$ws.Range("D199").Value = 44651
$ws.Range("I199").Value = '1a (cosecha)'
$ws.Range("J199").Value = 600
$ws.Range("K199").Value = 300
$ws.Range("L199").Value = 350
$ws.Range("M199").Value = 325
$ws.Range("P199").Value = 325

# Row 200
$ws.Range("D200").Value = 44453
$ws.Range("H200").Value = 'Camote'
$ws.Range("J200").Value = 600
$ws.Range("K200").Value = 600
$ws.Range("L200").Value = 650
$ws.Range("M200").Value = 625
$ws.Range("P200").Value = 625

# Row 201
$ws.Range("D201").Value = 44453
$ws.Range("H201").Value = 'Camote'
$ws.Range("J201").Value = 300
$ws.Range("K201").Value = 550
$ws.Range("L201").Value = 550
$ws.Range("M201").Value = 550
$ws.Range("P201").Value = 550

# Row 202
$ws.Range("D202").Value = 44526
$ws.Range("H202").Value = 'Paine'
$ws.Range("J202").Value = 400
$ws.Range("K202").Value = 140
$ws.Range("L202").Value = 150
$ws.Range("M202").Value = 145
$ws.Range("P202").Value = 145

# Row 203
$ws.Range("D203").Value = 44526
$ws.Range("H203").Value = 'Paine'
$ws.Range("J203").Value = 200
$ws.Range("K203").Value = 100
$ws.Range("L203").Value = 100
$ws.Range("M203").Value = 100
$ws.Range("P203").Value = 100

# Row 204
$ws.Range("D204").Value = 44363
$ws.Range("I204").Value = '1a (guarda)'
$ws.Range("K204").Value = 350
$ws.Range("L204").Value = 360
$ws.Range("M204").Value = 355
$ws.Range("P204").Value = 355

# Row 205
$ws.Range("D205").Value = 44363
$ws.Range("I205").Value = '2a (guarda)'
$ws.Range("K205").Value = 250
$ws.Range("L205").Value = 250
$ws.Range("M205").Value = 250
$ws.Range("P205").Value = 250

# Row 206
$ws.Range("D206").Value = 44251
$ws.Range("I206").Value = '1a nueva(o)'
$ws.Range("J206").Value = 600
$ws.Range("K206").Value = 250
$ws.Range("L206").Value = 300
$ws.Range("M206").Value = 275
$ws.Range("P206").Value = 275

# Row 207
$ws.Range("D207").Value = 44251
$ws.Range("I207").Value = '2a nueva(o)'
$ws.Range("J207").Value = 300
$ws.Range("K207").Value = 200
$ws.Range("L207").Value = 200
$ws.Range("M207").Value = 200
$ws.Range("P207").Value = 200

# Row 208
$ws.Range("D208").Value = 44617
$ws.Range("H208").Value = 'Camote'
$ws.Range("I208").Value = '1a (cosecha)'
$ws.Range("J208").Value = 400
$ws.Range("K208").Value = 500
$ws.Range("L208").Value = 500
$ws.Range("M208").Value = 500
$ws.Range("P208").Value = 500

# Row 209
$ws.Range("D209").Value = 44617
$ws.Range("H209").Value = 'Camote'
$ws.Range("I209").Value = '2a (cosecha)'
$ws.Range("J209").Value = 400
$ws.Range("K209").Value = 400
$ws.Range("L209").Value = 400
$ws.Range("M209").Value = 400
$ws.Range("P209").Value = 400

# Row 210
$ws.Range("D210").Value = 44169
$ws.Range("H210").Value = 'Paine'
$ws.Range("I210").Value = '1a nueva(o)'
$ws.Range("J210").Value = 500
$ws.Range("K210").Value = 800
$ws.Range("L210").Value = 800
$ws.Range("M210").Value = 800
$ws.Range("P210").Value = 800

# Row 211
$ws.Range("D211").Value = 44169
$ws.Range("H211").Value = 'Paine'
$ws.Range("I211").Value = '2a nueva(o)'
$ws.Range("K211").Value = 700
$ws.Range("L211").Value = 700
$ws.Range("M211").Value = 700
$ws.Range("P211").Value = 700

# Row 212
$ws.Range("A212").Value = 11
$ws.Range("B212").Value = 'Vega Monumental Concepción'
$ws.Range("C212").Value = 'Bíobío'
$ws.Range("D212").Value = 44474
$ws.Range("E212").Value = 8
$ws.Range("F212").Value = 100112045
$ws.Range("G212").Value = 'Zapallo'
$ws.Range("H212").Value = 'Camote'
$ws.Range("I212").Value = '1a (guarda)'
$ws.Range("J212").Value = 600
$ws.Range("K212").Value = 600
$ws.Range("L212").Value = 650
$ws.Range("M212").Value = 625
$ws.Range("N212").Value = '$/kilo (volumen en unidades)'
$ws.Range("O212").Value = 'Región de O''Higgins'
$ws.Range("P212").Value = 625
$ws.Range("Q212").Value = 1
$ws.Range("R212").Value = 'Hortaliza'

# Row 213
$ws.Range("A213").Value = 11
$ws.Range("B213").Value = 'Vega Monumental Concepción'
$ws.Range("C213").Value = 'Bíobío'
$ws.Range("D213").Value = 44474
$ws.Range("E213").Value = 8
$ws.Range("F213").Value = 100112045
$ws.Range("G213").Value = 'Zapallo'
$ws.Range("H213").Value = 'Camote'
$ws.Range("I213").Value = '2a (guarda)'
$ws.Range("J213").Value = 300
$ws.Range("K213").Value = 550
$ws.Range("L213").Value = 550
$ws.Range("M213").Value = 550
$ws.Range("N213").Value = '$/kilo (volumen en unidades)'
$ws.Range("O213").Value = 'Región de O''Higgins'
$ws.Range("P213").Value = 550
$ws.Range("Q213").Value = 1
$ws.Range("R213").Value = 'Hortaliza'
